# fix(super_admin): permitir administradores cargar copropietarios
#
# The "password" column (G) is removed from the user-upload template
# (users should no longer need to supply a password in the bulk-upload
# sheet), and the sample row's apartment_number is updated from 101 to
# 102. The active selection is also moved to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample apartment number (column E) from 101 to 102.
$ws.Range("E2").Value = 102

# Remove the whole "password" column (G1 header + G2 sample value).
# This shifts nothing else since it is the last column, and updates the
# sheet dimension / shared strings automatically.
$ws.Range("G1:G2").EntireColumn.Delete()

# Match the saved selection/active cell from the authored workbook.
[void]$ws.Range("E5").Select()
